$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dezembro")

$eurFormat = "#,##0.00 €; [Red]-#,##0.00 €"

# --- Row 7: the expense that used to live in row 8 ("Health") now occupies
# row 7; its date (04/12/2023) already matches what row 7 had, so only the
# amount/purpose/description/category need to change.
$ws.Range("B7").NumberFormat = $eurFormat
$ws.Range("B7").Value = -12
$ws.Range("C7").Value = "sdgdf"
$ws.Range("D7").Value = "gsfgbfg"
$ws.Range("E7").Value = "Health"

# --- Row 8: brand-new "Education" expense entry added in this revision.
# Force the date to be stored as plain text (matching the sheet's existing
# convention for the DATE column), not auto-converted to a date serial
# number by Excel's input parser.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "07/12/2023"

$ws.Range("B8").NumberFormat = $eurFormat
$ws.Range("B8").Value = -20
$ws.Range("C8").Value = "kdsjfa"
$ws.Range("D8").Value = "jklsndnvs"
$ws.Range("E8").Value = "Education"

# --- Row 9: the expense that used to live in row 7 ("af" / 231) now occupies
# row 9, revised into the new "Other" entry. Same date as row 7/8.
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "04/12/2023"

$ws.Range("B9").NumberFormat = $eurFormat
$ws.Range("B9").Value = 100
$ws.Range("C9").Value = "pai"
$ws.Range("D9").Value = "ola"
$ws.Range("E9").Value = "Other"

# Drop the temporary "Text" number format now that the literal date strings
# are safely stored, so A8/A9 end up with the plain default style - just
# like every other DATE cell in this column (A4:A8 originally had none).
$ws.Range("A8:A9").ClearFormats()
